# Updates the crypto price/volume(1h) snapshot on Sheet1 to the new data
# pull (per the "Updated cryptos list" GitHub Actions commit).
# Column D ("Price") holds text-formatted numbers (e.g. "226.39" or the
# dotted-thousands form "36.290.89"); a leading apostrophe forces Excel to
# store it as literal text instead of normalizing/coercing it to a number.
# Column E ("Volume(1h)") values already contain non-numeric characters
# (% and padding spaces) so Excel keeps them as text without any trick.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''36.290.89'
$ws.Range("E2").Value = '  -3.34%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''1.951.00'
$ws.Range("E3").Value = '  -3.65%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.02%  '

# Row 5: BNB
$ws.Range("D5").Value = '''226.39'
$ws.Range("E5").Value = '  -11.19%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.589'
$ws.Range("E6").Value = '  -5.01%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.04%  '

# Row 8: Solana
$ws.Range("D8").Value = '''52.50'
$ws.Range("E8").Value = '  -7.53%  '

# Row 9: Cardano
$ws.Range("E9").Value = '  -5.42%  '

# Row 10: OKB
$ws.Range("D10").Value = '''56.61'
$ws.Range("E10").Value = '  -1.33%  '

# Row 11: Dogecoin
$ws.Range("E11").Value = '  -7.25%  '

# Row 12: TRON
$ws.Range("D12").Value = '''0.0966'
$ws.Range("E12").Value = '  -4.56%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '''2.235.90'
$ws.Range("E13").Value = '  -3.82%  '

# Row 14: Chainlink
$ws.Range("D14").Value = '''13.61'
$ws.Range("E14").Value = '  -6.19%  '

# Row 15: Avalanche
$ws.Range("D15").Value = '''19.24'
$ws.Range("E15").Value = '  -8.59%  '

# Row 16: Polygon
$ws.Range("D16").Value = '''0.734'
$ws.Range("E16").Value = '  -9.98%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '''1.952.61'
$ws.Range("E17").Value = '  -3.59%  '

# Row 18: Polkadot
$ws.Range("D18").Value = '''4.92'
$ws.Range("E18").Value = '  -7.93%  '

# Row 19: WrappedBTC
$ws.Range("D19").Value = '''36.218.77'
$ws.Range("E19").Value = '  -3.36%  '

# Row 20: Litecoin
$ws.Range("D20").Value = '''66.41'
$ws.Range("E20").Value = '  -4.52%  '

# Row 21: ShibaInu
$ws.Range("D21").Value = '''0.0₃0781'
$ws.Range("E21").Value = '  -7.82%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  -5.52%  '

# Row 23: BitcoinCash
$ws.Range("D23").Value = '''217.96'
$ws.Range("E23").Value = '  -4.59%  '

# Row 24: Dai
$ws.Range("E24").Value = '  +0.09%  '

# Row 25: Toncoin
$ws.Range("D25").Value = '''2.33'
$ws.Range("E25").Value = '  +0.20%  '

# Row 26: PancakeSwap
$ws.Range("E26").Value = '  -12.73%  '

# Row 27: Monero
$ws.Range("D27").Value = '''159.38'
$ws.Range("E27").Value = '  -2.48%  '

# Row 28: Cosmos
$ws.Range("D28").Value = '''8.30'
$ws.Range("E28").Value = '  -8.19%  '

# Row 29: EthereumClassic
$ws.Range("D29").Value = '''18.50'
$ws.Range("E29").Value = '  -6.64%  '

# Row 30: Kaspa
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '''1.27'
$ws.Range("E30").Value = '  -7.05%  '

# Row 31: ImmutableX
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '''0.117'
$ws.Range("E31").Value = '  -11.45%  '

# Row 32: Stellar
$ws.Range("D32").Value = '''0.115'
$ws.Range("E32").Value = '  -4.83%  '

# Row 33: Filecoin
$ws.Range("E33").Value = '  -9.44%  '

# Row 34: Hedera
$ws.Range("D34").Value = '''0.0593'
$ws.Range("E34").Value = '  -11.12%  '

# Row 35: InternetComputer(DFINITY)
$ws.Range("E35").Value = '  -9.49%  '

# Row 36: BinanceUSD
$ws.Range("E36").Value = '  +0.10%  '

# Row 37: LidoDAOToken
$ws.Range("E37").Value = '  -8.70%  '

# Row 38: WEMIXToken
$ws.Range("E38").Value = '  -2.15%  '

# Row 39: RenderToken
$ws.Range("D39").Value = '''3.09'
$ws.Range("E39").Value = '  -8.59%  '

# Row 40: HuobiToken
$ws.Range("D40").Value = '''2.99'
$ws.Range("E40").Value = '  -1.25%  '

# Row 41: THORChain
$ws.Range("D41").Value = '''4.94'
$ws.Range("E41").Value = '  -7.78%  '

# Row 42: Maker
$ws.Range("D42").Value = '''1.392.60'
$ws.Range("E42").Value = '  -1.13%  '

# Row 43: VeChain
$ws.Range("E43").Value = '  -9.14%  '

# Row 44: Cronos
$ws.Range("D44").Value = '''0.0853'
$ws.Range("E44").Value = '  -11.69%  '

# Row 45: TrustWalletToken
$ws.Range("E45").Value = '  -11.41%  '

# Row 46: Aave
$ws.Range("D46").Value = '''85.22'
$ws.Range("E46").Value = '  -6.03%  '

# Row 47: MXToken
$ws.Range("D47").Value = '''2.84'
$ws.Range("E47").Value = '  -1.00%  '

# Row 48: ARBITRUM
$ws.Range("E48").Value = '  -7.02%  '

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''14.41'
$ws.Range("E49").Value = '  -9.76%  '

# Row 50: FraxShare
$ws.Range("E50").Value = '  -9.37%  '

# Row 51: RocketPoolETH
$ws.Range("D51").Value = '''2.129.95'
$ws.Range("E51").Value = '  -3.93%  '
